# Adds a new "Sheet8" at the end of the workbook containing a small
# Date / Float / Int / Time reference table (the xlrd date/time
# interpretation test fixture referenced in the commit message), and makes
# it the active sheet/tab (mirrors Excel's "tabSelected" moving from the
# previously-last sheet to the freshly added one).

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current last sheet so it lands at
# the end of the tab strip (Worksheets.Add() with no "After"/"Before"
# argument would insert it *before* the active sheet instead).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Sheet8"

# Row 1: Date
$ws.Range("A1").Value = "Date:"
$ws.Range("B1").Value = 43070
$ws.Range("B1").NumberFormat = "d-mmm-yy"

# Row 2: Float
$ws.Range("A2").Value = "Float:"
$ws.Range("B2").Value = 1.12

# Row 3: Int
$ws.Range("A3").Value = "Int:"
$ws.Range("B3").Value = 2

# Row 4: Time
$ws.Range("A4").Value = "Time:"
$ws.Range("B4").Value = 0.10416666666666667
$ws.Range("B4").NumberFormat = "h:mm"

# Column B is best-fit to hold the longest formatted value ("1-Dec-17").
$ws.Columns.Item(2).ColumnWidth = 8.8

[void]$ws.Range("A4:B4").Select()
